# Append the latest daily mod-count reading to the ModCounts sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModCounts")

# Find the last populated row (mirrors the existing data block starting at row 2)
# and write the new record directly below it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Column A holds plain date-label text (e.g. "2026/01/28"), not a real date
# serial, so force the cell to Text first - otherwise Excel would silently
# reinterpret "2026/01/29" as a date value.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026/01/29"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1164

# Match the formatting of the preceding row (centered alignment etc.)
# instead of leaving the text-forced cell with its own one-off style.
$ws.Range("A${lastRow}:C${lastRow}").Copy()
$ws.Range("A${newRow}:C${newRow}").PasteSpecial(-4122)
